$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the first 4 data rows (rows 2-5), shifting the remaining data rows up
$ws.Range("A2:C5").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftUp)

# Append 14 new rows of data after the shifted data (rows 18-31)
$ws.Cells.Item(18, 1).Value = 17.96537399291992
$ws.Cells.Item(18, 2).Value = -41.82672500610352
$ws.Cells.Item(18, 3).Value = 13.37196350097656
$ws.Cells.Item(19, 1).Value = 16.90632247924805
$ws.Cells.Item(19, 2).Value = 11.39815711975098
$ws.Cells.Item(19, 3).Value = -9.546463966369627
$ws.Cells.Item(20, 1).Value = -13.79934120178223
$ws.Cells.Item(20, 2).Value = -5.43239688873291
$ws.Cells.Item(20, 3).Value = -14.25337219238281
$ws.Cells.Item(21, 1).Value = -7.913525104522705
$ws.Cells.Item(21, 2).Value = -43.03920745849609
$ws.Cells.Item(21, 3).Value = 5.675649166107178
$ws.Cells.Item(22, 1).Value = -19.80719947814941
$ws.Cells.Item(22, 2).Value = 17.96865081787109
$ws.Cells.Item(22, 3).Value = 0.9991121292114258
$ws.Cells.Item(23, 1).Value = 5.307936668395996
$ws.Cells.Item(23, 2).Value = -36.55269622802734
$ws.Cells.Item(23, 3).Value = 44.45186233520508
$ws.Cells.Item(24, 1).Value = -72.68199157714844
$ws.Cells.Item(24, 2).Value = -20.47943878173828
$ws.Cells.Item(24, 3).Value = -39.62965774536133
$ws.Cells.Item(25, 1).Value = -8.546328544616699
$ws.Cells.Item(25, 2).Value = 0.6573230028152466
$ws.Cells.Item(25, 3).Value = -10.66314029693604
$ws.Cells.Item(26, 1).Value = -15.78201103210449
$ws.Cells.Item(26, 2).Value = -11.37719249725342
$ws.Cells.Item(26, 3).Value = -7.033545017242432
$ws.Cells.Item(27, 1).Value = -11.98593044281006
$ws.Cells.Item(27, 2).Value = 2.52801513671875
$ws.Cells.Item(27, 3).Value = 9.990962028503418
$ws.Cells.Item(28, 1).Value = -9.182360649108888
$ws.Cells.Item(28, 2).Value = 19.98668098449707
$ws.Cells.Item(28, 3).Value = 6.591425895690918
$ws.Cells.Item(29, 1).Value = 14.62302780151367
$ws.Cells.Item(29, 2).Value = -34.03761672973633
$ws.Cells.Item(29, 3).Value = 14.29042530059814
$ws.Cells.Item(30, 1).Value = -34.80189895629883
$ws.Cells.Item(30, 2).Value = 14.59354972839356
$ws.Cells.Item(30, 3).Value = -21.27434158325196
$ws.Cells.Item(31, 1).Value = -21.16343307495117
$ws.Cells.Item(31, 2).Value = -7.425243377685547
$ws.Cells.Item(31, 3).Value = -10.02677822113037
